$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77 (shifts existing rows 77-192 down to 78-193,
# carrying row formatting - e.g. the date style on column D - down with them).
$ws.Rows(77).Insert()

# Populate the newly inserted row 77 with the new weekly record.
$ws.Range("A77").Value = 11
$ws.Range("B77").Value = "Vega Monumental Concepción"
$ws.Range("C77").Value = "Bíobío"
$ws.Range("D77").Value = 45174
$ws.Range("E77").Value = 8
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108002
$ws.Range("J77").Value = "Mango"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 200
$ws.Range("N77").Value = 9000
$ws.Range("O77").Value = 10000
$ws.Range("P77").Value = 9500
$ws.Range("Q77").Value = "$/bandeja 4 kilos"
$ws.Range("R77").Value = "Brasil"
$ws.Range("S77").Value = 2375
$ws.Range("T77").Value = 4
